$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Purchase Request field (merged C7:E7): Bacolod -> SITE
$ws.Range("C7").Value = "SITE"

# Department field (merged I7:K7): IT -> ITSite
$ws.Range("I7").Value = "ITSite"

# Date Prepared (merged C8:E8): 2020-02-22 -> 2020-06-02
$ws.Range("C8").Value = "2020-06-02"

# Dept. Code dropdown (merged I8:K8): ITB -> ITS (IT Department - SITE)
$ws.Range("I8").Value = "ITS"

# Requestor (merged I9:K9): Stephine David -> Jushkyle Jambongana
$ws.Range("I9").Value = "Jushkyle Jambongana"

# Row 14: Qty 10 -> 5
$ws.Range("B14").Value = 5

# Row 15: Item No 2 -> 1, Qty 1 -> 5
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = 5

# Move the active selection to L14 (matches final view state)
$ws.Range("L14").Select()
